$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 90,3
$data[0,0] = 0; $data[0,1] = 0.09; $data[0,2] = 0.91
$data[1,0] = 0; $data[1,1] = 1; $data[1,2] = 0
$data[2,0] = 0.74; $data[2,1] = 0.13; $data[2,2] = 0.13
$data[3,0] = 1; $data[3,1] = 0; $data[3,2] = 0
$data[4,0] = 0; $data[4,1] = 1; $data[4,2] = 0
$data[5,0] = 0; $data[5,1] = 1; $data[5,2] = 0
$data[6,0] = 0; $data[6,1] = 0; $data[6,2] = 1
$data[7,0] = 0; $data[7,1] = 1; $data[7,2] = 0
$data[8,0] = 0; $data[8,1] = 1; $data[8,2] = 0
$data[9,0] = 0.04; $data[9,1] = 0.08; $data[9,2] = 0.89
$data[10,0] = 0; $data[10,1] = 1; $data[10,2] = 0
$data[11,0] = 0; $data[11,1] = 1; $data[11,2] = 0
$data[12,0] = 0; $data[12,1] = 1; $data[12,2] = 0
$data[13,0] = 1; $data[13,1] = 0; $data[13,2] = 0
$data[14,0] = 0; $data[14,1] = 0; $data[14,2] = 1
$data[15,0] = 0; $data[15,1] = 0; $data[15,2] = 1
$data[16,0] = 1; $data[16,1] = 0; $data[16,2] = 0
$data[17,0] = 1; $data[17,1] = 0; $data[17,2] = 0
$data[18,0] = 0; $data[18,1] = 1; $data[18,2] = 0
$data[19,0] = 1; $data[19,1] = 0; $data[19,2] = 0
$data[20,0] = 0; $data[20,1] = 0; $data[20,2] = 1
$data[21,0] = 0; $data[21,1] = 0.99; $data[21,2] = 0.01
$data[22,0] = 0; $data[22,1] = 1; $data[22,2] = 0
$data[23,0] = 1; $data[23,1] = 0; $data[23,2] = 0
$data[24,0] = 1; $data[24,1] = 0; $data[24,2] = 0
$data[25,0] = 1; $data[25,1] = 0; $data[25,2] = 0
$data[26,0] = 0; $data[26,1] = 0; $data[26,2] = 1
$data[27,0] = 0; $data[27,1] = 0; $data[27,2] = 1
$data[28,0] = 0; $data[28,1] = 1; $data[28,2] = 0
$data[29,0] = 0; $data[29,1] = 1; $data[29,2] = 0
$data[30,0] = 0; $data[30,1] = 1; $data[30,2] = 0
$data[31,0] = 1; $data[31,1] = 0; $data[31,2] = 0
$data[32,0] = 0; $data[32,1] = 1; $data[32,2] = 0
$data[33,0] = 1; $data[33,1] = 0; $data[33,2] = 0
$data[34,0] = 0; $data[34,1] = 0; $data[34,2] = 1
$data[35,0] = 1; $data[35,1] = 0; $data[35,2] = 0
$data[36,0] = 0; $data[36,1] = 1; $data[36,2] = 0
$data[37,0] = 0; $data[37,1] = 0.98; $data[37,2] = 0.02
$data[38,0] = 1; $data[38,1] = 0; $data[38,2] = 0
$data[39,0] = 0; $data[39,1] = 1; $data[39,2] = 0
$data[40,0] = 0.95; $data[40,1] = 0.02; $data[40,2] = 0.02
$data[41,0] = 0; $data[41,1] = 1; $data[41,2] = 0
$data[42,0] = 1; $data[42,1] = 0; $data[42,2] = 0
$data[43,0] = 0; $data[43,1] = 0; $data[43,2] = 1
$data[44,0] = 1; $data[44,1] = 0; $data[44,2] = 0
$data[45,0] = 0; $data[45,1] = 1; $data[45,2] = 0
$data[46,0] = 0; $data[46,1] = 1; $data[46,2] = 0
$data[47,0] = 0; $data[47,1] = 1; $data[47,2] = 0
$data[48,0] = 0; $data[48,1] = 1; $data[48,2] = 0
$data[49,0] = 0; $data[49,1] = 0; $data[49,2] = 1
$data[50,0] = 0; $data[50,1] = 0; $data[50,2] = 1
$data[51,0] = 0; $data[51,1] = 0; $data[51,2] = 1
$data[52,0] = 0.12; $data[52,1] = 0.79; $data[52,2] = 0.09
$data[53,0] = 0; $data[53,1] = 0; $data[53,2] = 1
$data[54,0] = 0; $data[54,1] = 0; $data[54,2] = 1
$data[55,0] = 1; $data[55,1] = 0; $data[55,2] = 0
$data[56,0] = 0; $data[56,1] = 0; $data[56,2] = 1
$data[57,0] = 0; $data[57,1] = 0; $data[57,2] = 1
$data[58,0] = 0; $data[58,1] = 0.98; $data[58,2] = 0.02
$data[59,0] = 1; $data[59,1] = 0; $data[59,2] = 0
$data[60,0] = 0; $data[60,1] = 0.99; $data[60,2] = 0.01
$data[61,0] = 1; $data[61,1] = 0; $data[61,2] = 0
$data[62,0] = 1; $data[62,1] = 0; $data[62,2] = 0
$data[63,0] = 0; $data[63,1] = 1; $data[63,2] = 0
$data[64,0] = 0; $data[64,1] = 0; $data[64,2] = 1
$data[65,0] = 0; $data[65,1] = 0.99; $data[65,2] = 0.01
$data[66,0] = 0.33; $data[66,1] = 0.66; $data[66,2] = 0.01
$data[67,0] = 0; $data[67,1] = 1; $data[67,2] = 0
$data[68,0] = 0; $data[68,1] = 0; $data[68,2] = 1
$data[69,0] = 0; $data[69,1] = 0; $data[69,2] = 1
$data[70,0] = 1; $data[70,1] = 0; $data[70,2] = 0
$data[71,0] = 0; $data[71,1] = 0.99; $data[71,2] = 0.01
$data[72,0] = 0; $data[72,1] = 1; $data[72,2] = 0
$data[73,0] = 0; $data[73,1] = 0; $data[73,2] = 1
$data[74,0] = 1; $data[74,1] = 0; $data[74,2] = 0
$data[75,0] = 0; $data[75,1] = 0; $data[75,2] = 1
$data[76,0] = 0; $data[76,1] = 0; $data[76,2] = 1
$data[77,0] = 1; $data[77,1] = 0; $data[77,2] = 0
$data[78,0] = 0; $data[78,1] = 0; $data[78,2] = 1
$data[79,0] = 0; $data[79,1] = 0.72; $data[79,2] = 0.28
$data[80,0] = 0; $data[80,1] = 1; $data[80,2] = 0
$data[81,0] = 0; $data[81,1] = 0; $data[81,2] = 1
$data[82,0] = 0; $data[82,1] = 1; $data[82,2] = 0
$data[83,0] = 1; $data[83,1] = 0; $data[83,2] = 0
$data[84,0] = 0; $data[84,1] = 0; $data[84,2] = 1
$data[85,0] = 0.32; $data[85,1] = 0.13; $data[85,2] = 0.55
$data[86,0] = 1; $data[86,1] = 0; $data[86,2] = 0
$data[87,0] = 1; $data[87,1] = 0; $data[87,2] = 0
$data[88,0] = 0; $data[88,1] = 0; $data[88,2] = 1
$data[89,0] = 0; $data[89,1] = 0; $data[89,2] = 0

$ws.Range("B2:D91").Value = $data
